# ---------------------------------------------------------------------------
# Add "arrayCode" and "Intro" worksheets with their python-exercise content,
# matching the target commit ("Arraypage, DSintro, Basepage, datasheet").
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

$code25 = @'
def search(list,element):
    if element in list:
        return("Element Found")
    else:
        return("Not Found")
li=[9, 2, 34, 2, 5, 9, 16]
element= 9
print("Check if",element,"is in",li)
print(search(li,element))
'@

$code26 = @'
def findMaxConsecutiveOnes(nums):
    count = 0
    result = 0
    for i in range(0, len(nums)):
        if (nums[i] == 0):
            count = 0
        else:
            count+= 1
            result = max(result, count)
    return result
nums  = [1,0,1,1,0,1]
n = len(nums)
print(findMaxConsecutiveOnes(nums))
'@

$code27 = @'
def findNumbers(nums):
	c=0
	for i in nums:
		j=str(i)
		x=len(j)
		if x%2==0:
				c=c+1
        print (c)
	return c
findNumbers([12,345,2,6,7896])
'@

$code28 = @'
def findNumbers(nums):
	c=0
	for i in nums:
		j=str(i)
		x=len(j)
		if x%2==0:
				c=c+1
	return c
findNumbers([12,345,2,6,7896])
'@

$code29 = @'
def sortedSquares(nums):
	squares_list = []
 	for i in range(0, len(nums)):
     		square = nums[i] * nums[i];
     		squares_list.append(square)
     		sorted_squares_list = sorted(squares_list)
   	print sorted_squares_list;	
	return sorted_squares_list; 		
sortedSquares([-7,-3,2,3,11])
'@

$code31 = @'
def search(list,element):
    if element in list:
        return("Element Found")
    else:
        return("Not Found")
li=[9, 2, 34, 2, 5, 9, 16]
element= 9
print("Check if",element,"is in",li)
print(search(li,element))
'@


# ---------------------------------------------------------------------------
# 1. Add the two new worksheets at the end, in order: arrayCode, then Intro
# ---------------------------------------------------------------------------
$wsArray = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsArray.Name = "arrayCode"

$wsIntro = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsIntro.Name = "Intro"

# ---------------------------------------------------------------------------
# 2. arrayCode sheet content
# ---------------------------------------------------------------------------
$wsArray.Range("A1").Value = "arraycode"
$wsArray.Range("B1").Value = "Result"

$wsArray.Range("A2").Value = 'print("Hello World")'
$wsArray.Range("B2").Value = "Hello World"

$wsArray.Range("A3").Value = "Hello World"
$wsArray.Range("B3").Value = "SyntaxError: bad input on line 1"
$wsArray.Range("B3").Font.Name = "Courier New"
$wsArray.Range("B3").Font.Size = 10
$wsArray.Range("B3").Font.Color = 0
$wsArray.Range("B3").Font.Family = 3
$wsArray.Range("B3").VerticalAlignment = -4108

$wsArray.Range("A4").Value = $code31
$wsArray.Range("B4").Value = "Element Found"

$wsArray.Range("A5").Value = $code25
$wsArray.Range("B5").Value = "submission success"

$wsArray.Range("A6").Value = $code26
$wsArray.Range("B6").Value = "2"

$wsArray.Range("A7").Value = $code26
$wsArray.Range("B7").Value = "submission success"

$wsArray.Range("A8").Value = $code27
$wsArray.Range("B8").Value = "2"

$wsArray.Range("A9").Value = $code28
$wsArray.Range("B9").Value = "submission success"

$wsArray.Range("B10").Value = "[4, 9, 9, 49, 121]"

$wsArray.Range("A11").Value = $code29
$wsArray.Range("B11").Value = "No tests were collected"

# Number format + font for the "code" column (A4:A11) and "result" column (B5:B11)
$codeRange = $wsArray.Range("A4:A11")
$codeRange.Font.Name = "Calibri"
$codeRange.Font.Size = 11
$codeRange.Font.Color = 0
$codeRange.WrapText = $true
$codeRange.HorizontalAlignment = -4131
$codeRange.Interior.Color = 16777215
$codeRange.Interior.PatternColor = 16777215

$resultCodeRange = $wsArray.Range("B5:B11")
$resultCodeRange.NumberFormat = "@"
$resultCodeRange.Font.Name = "Calibri"
$resultCodeRange.Font.Size = 11

$wsArray.Range("B4").Font.Name = "Calibri"
$wsArray.Range("B4").Font.Size = 11

# Column widths (character widths tuned so the exported pixel-quantised
# width lands on the target stored width)
$wsArray.Columns.Item(1).ColumnWidth = 42.17
$wsArray.Columns.Item(2).ColumnWidth = 37.67

# Row heights for the wrapped code rows
$wsArray.Rows.Item(4).RowHeight = 165
$wsArray.Rows.Item(5).RowHeight = 135
$wsArray.Rows.Item(6).RowHeight = 270
$wsArray.Rows.Item(7).RowHeight = 270
$wsArray.Rows.Item(8).RowHeight = 150
$wsArray.Rows.Item(9).RowHeight = 150
$wsArray.Rows.Item(11).RowHeight = 150

$wsArray.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 3. Intro sheet content
# ---------------------------------------------------------------------------
$wsIntro.Range("A1").Value = "introcode"
$wsIntro.Range("B1").Value = "Result"

$wsIntro.Range("A2").Value = 'print("Hello World")'
$wsIntro.Range("B2").Value = "Hello World"

$wsIntro.Range("A3").Value = "Hello World"
$wsIntro.Range("B3").Value = "NameError:"

$wsIntro.Columns.Item(1).ColumnWidth = 18.505
$wsIntro.Columns.Item(2).ColumnWidth = 28.335

# ---------------------------------------------------------------------------
# 4. Selection / active-sheet bookkeeping (mirrors the diff's sheetViews)
# ---------------------------------------------------------------------------
$wsIntro.Range("B2").Select() | Out-Null
$wsArray.Activate()
$wsArray.Range("A4").Select() | Out-Null
